$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.077.95'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.653.46'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5251'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2600'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06335'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.36'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07796'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.500'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("D13").Value = '1.649.12'
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5473'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '0.0₅8238'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").Value = '26.107.67'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.572'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '190.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.024'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '141.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1230'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.233'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05865'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.532'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.260'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.579'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9486'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.412'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.779'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5720'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01615'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.778'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8454'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.82%  '
$ws.Range("D43").Value = '1.026.40'
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").Value = '1.799.42'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4314'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05150'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.843'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.467'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09658'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.54%  '
